$wb = $excel.ActiveWorkbook

# ---- Sheet "Risks" (sheet index 2): append rows 19-38 ----
$risks = $wb.Worksheets.Item("Risks")

$risksData = @(
    @("e6f30fb4-a987-4631-934b-09b1829960fd", "2025-05-21T16:18:46.857457", 0.25, 0.5, 0, "MEDIUM", "COMPLETED", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "Enable drift monitoring; schedule periodic retraining", "article_9", "COMPLETED", "2025-05-21T16:18:46.857457"),
    @("b69eac1e-8827-4459-92c2-78d08474f636", "2025-05-21T17:33:25.088759", 0.25, 0.5, 0, "MEDIUM", "COMPLETED", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "Enable drift monitoring; schedule periodic retraining", "article_9", "COMPLETED", "2025-05-21T17:33:25.088759"),
    @("06b5dd9e-f246-48a6-9cac-86ff38f32d76", "2025-05-21T17:46:22.866806", 0.25, 0.5, 0, "MEDIUM", "COMPLETED", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "Enable drift monitoring; schedule periodic retraining", "article_9", "COMPLETED", "2025-05-21T17:46:22.866806"),
    @("23fff188-c858-4c00-925d-213db01ea520", "2025-05-21T18:03:11.593143", 0.25, 0.5, 0, "MEDIUM", "COMPLETED", "Model accuracy below 0.75", "Collect more data; tune hyper-parameters", "article_9", "COMPLETED", "2025-05-21T18:03:11.593143"),
    @("23fff188-c858-4c00-925d-213db01ea520", "2025-05-21T18:03:11.593143", 0.25, 0.5, 0, "MEDIUM", "COMPLETED", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "Enable drift monitoring; schedule periodic retraining", "article_9", "COMPLETED", "2025-05-21T18:03:11.593143"),
    @("82422ddc-71ec-4371-9842-95dc1db0b38f", "2025-05-21T18:08:03.707097", 0.65, 0.5, 0.8, "HIGH", "PENDING", "Unfair bias against protected demographic groups", "Re-sample training data; add fairness constraints or post-processing techniques", "article_10", "PENDING", "2025-05-21T18:08:03.707097"),
    @("82422ddc-71ec-4371-9842-95dc1db0b38f", "2025-05-21T18:08:03.707097", 0.65, 0.5, 0.8, "MEDIUM", "PENDING", "Model accuracy below 0.75", "Collect more data; tune hyper-parameters", "article_9", "PENDING", "2025-05-21T18:08:03.707097"),
    @("82422ddc-71ec-4371-9842-95dc1db0b38f", "2025-05-21T18:08:03.707097", 0.65, 0.5, 0.8, "MEDIUM", "PENDING", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "Enable drift monitoring; schedule periodic retraining", "article_9", "PENDING", "2025-05-21T18:08:03.707097"),
    @("35b3012a-e79a-46bc-b468-2fda279cfe8c", "2025-05-21T18:12:25.589810", 0.65, 0.5, 0.8, "HIGH", "PENDING", "Unfair bias against protected demographic groups", "Re-sample training data; add fairness constraints or post-processing techniques", "article_10", "PENDING", "2025-05-21T18:12:25.589810"),
    @("35b3012a-e79a-46bc-b468-2fda279cfe8c", "2025-05-21T18:12:25.589810", 0.65, 0.5, 0.8, "MEDIUM", "PENDING", "Model accuracy below 0.75", "Collect more data; tune hyper-parameters", "article_9", "PENDING", "2025-05-21T18:12:25.589810"),
    @("35b3012a-e79a-46bc-b468-2fda279cfe8c", "2025-05-21T18:12:25.589810", 0.65, 0.5, 0.8, "MEDIUM", "PENDING", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "Enable drift monitoring; schedule periodic retraining", "article_9", "PENDING", "2025-05-21T18:12:25.589810"),
    @("b4421b8e-7d46-4cc2-9b4f-effdb9f3e8a4", "2025-05-21T18:13:46.120087", 0.65, 0.5, 0.8, "HIGH", "PENDING", "Unfair bias against protected demographic groups", "Re-sample training data; add fairness constraints or post-processing techniques", "article_10", "PENDING", "2025-05-21T18:13:46.120087"),
    @("b4421b8e-7d46-4cc2-9b4f-effdb9f3e8a4", "2025-05-21T18:13:46.120087", 0.65, 0.5, 0.8, "MEDIUM", "PENDING", "Model accuracy below 0.75", "Collect more data; tune hyper-parameters", "article_9", "PENDING", "2025-05-21T18:13:46.120087"),
    @("b4421b8e-7d46-4cc2-9b4f-effdb9f3e8a4", "2025-05-21T18:13:46.120087", 0.65, 0.5, 0.8, "MEDIUM", "PENDING", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "Enable drift monitoring; schedule periodic retraining", "article_9", "PENDING", "2025-05-21T18:13:46.120087"),
    @("ea9382b8-ad1d-41a8-a10b-9b28715e0f77", "2025-05-21T19:06:22.324942", 0.65, 0.5, 0.8, "HIGH", "PENDING", "Unfair bias against protected demographic groups", "Re-sample training data; add fairness constraints or post-processing techniques", "article_10", "PENDING", "2025-05-21T19:06:22.324942"),
    @("ea9382b8-ad1d-41a8-a10b-9b28715e0f77", "2025-05-21T19:06:22.324942", 0.65, 0.5, 0.8, "MEDIUM", "PENDING", "Model accuracy below 0.75", "Collect more data; tune hyper-parameters", "article_9", "PENDING", "2025-05-21T19:06:22.324942"),
    @("ea9382b8-ad1d-41a8-a10b-9b28715e0f77", "2025-05-21T19:06:22.324942", 0.65, 0.5, 0.8, "MEDIUM", "PENDING", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "Enable drift monitoring; schedule periodic retraining", "article_9", "PENDING", "2025-05-21T19:06:22.324942"),
    @("e244f4ce-832c-4618-addc-984857f50653", "2025-05-21T20:19:38.154359", 0.65, 0.5, 0.8, "HIGH", "PENDING", "Unfair bias against protected demographic groups", "Re-sample training data; add fairness constraints or post-processing techniques", "article_10", "PENDING", "2025-05-21T20:19:38.154359"),
    @("e244f4ce-832c-4618-addc-984857f50653", "2025-05-21T20:19:38.154359", 0.65, 0.5, 0.8, "MEDIUM", "PENDING", "Model accuracy below 0.75", "Collect more data; tune hyper-parameters", "article_9", "PENDING", "2025-05-21T20:19:38.154359"),
    @("e244f4ce-832c-4618-addc-984857f50653", "2025-05-21T20:19:38.154359", 0.65, 0.5, 0.8, "MEDIUM", "PENDING", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "Enable drift monitoring; schedule periodic retraining", "article_9", "PENDING", "2025-05-21T20:19:38.154359"),
)

$rowNum = 19
foreach ($row in $risksData) {
    $risks.Cells.Item($rowNum, 1).Value = $row[0]
    $risks.Cells.Item($rowNum, 2).Value = $row[1]
    $risks.Cells.Item($rowNum, 3).Value = $row[2]
    $risks.Cells.Item($rowNum, 4).Value = $row[3]
    $risks.Cells.Item($rowNum, 5).Value = $row[4]
    $risks.Cells.Item($rowNum, 6).Value = $row[5]
    $risks.Cells.Item($rowNum, 7).Value = $row[6]
    $risks.Cells.Item($rowNum, 8).Value = $row[7]
    $risks.Cells.Item($rowNum, 9).Value = $row[8]
    $risks.Cells.Item($rowNum, 10).Value = $row[9]
    $risks.Cells.Item($rowNum, 11).Value = $row[10]
    $risks.Cells.Item($rowNum, 12).Value = $row[11]
    $rowNum = $rowNum + 1
}

# ---- Sheet "HazardDetails" (sheet index 3): append rows 19-38 ----
$hazards = $wb.Worksheets.Item("HazardDetails")

$hazardsData = @(
    @("e6f30fb4-a987-4631-934b-09b1829960fd", "2025-05-21T16:18:46.857457", 0.25, "DRIFT_VULNERABILITY", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "MEDIUM", "Enable drift monitoring; schedule periodic retraining", $null, "article_9"),
    @("b69eac1e-8827-4459-92c2-78d08474f636", "2025-05-21T17:33:25.088759", 0.25, "DRIFT_VULNERABILITY", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "MEDIUM", "Enable drift monitoring; schedule periodic retraining", $null, "article_9"),
    @("06b5dd9e-f246-48a6-9cac-86ff38f32d76", "2025-05-21T17:46:22.866806", 0.25, "DRIFT_VULNERABILITY", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "MEDIUM", "Enable drift monitoring; schedule periodic retraining", $null, "article_9"),
    @("23fff188-c858-4c00-925d-213db01ea520", "2025-05-21T18:03:11.593143", 0.25, "LOW_ACCURACY", "Model accuracy below 0.75", "MEDIUM", "Collect more data; tune hyper-parameters", $null, "article_9"),
    @("23fff188-c858-4c00-925d-213db01ea520", "2025-05-21T18:03:11.593143", 0.25, "DRIFT_VULNERABILITY", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "MEDIUM", "Enable drift monitoring; schedule periodic retraining", $null, "article_9"),
    @("82422ddc-71ec-4371-9842-95dc1db0b38f", "2025-05-21T18:08:03.707097", 0.65, "BIAS_PROTECTED_GROUPS", "Unfair bias against protected demographic groups", "HIGH", "Re-sample training data; add fairness constraints or post-processing techniques", "num__AGE_YEARS: 1.000 disparity`nNAME_EDUCATION_TYPE: 0.410 disparity`nNAME_HOUSING_TYPE: 0.202 disparity`n", "article_10"),
    @("82422ddc-71ec-4371-9842-95dc1db0b38f", "2025-05-21T18:08:03.707097", 0.65, "LOW_ACCURACY", "Model accuracy below 0.75", "MEDIUM", "Collect more data; tune hyper-parameters", $null, "article_9"),
    @("82422ddc-71ec-4371-9842-95dc1db0b38f", "2025-05-21T18:08:03.707097", 0.65, "DRIFT_VULNERABILITY", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "MEDIUM", "Enable drift monitoring; schedule periodic retraining", $null, "article_9"),
    @("35b3012a-e79a-46bc-b468-2fda279cfe8c", "2025-05-21T18:12:25.589810", 0.65, "BIAS_PROTECTED_GROUPS", "Unfair bias against protected demographic groups", "HIGH", "Re-sample training data; add fairness constraints or post-processing techniques", "num__AGE_YEARS: 1.000 disparity`nNAME_EDUCATION_TYPE: 0.410 disparity`nNAME_HOUSING_TYPE: 0.202 disparity`n", "article_10"),
    @("35b3012a-e79a-46bc-b468-2fda279cfe8c", "2025-05-21T18:12:25.589810", 0.65, "LOW_ACCURACY", "Model accuracy below 0.75", "MEDIUM", "Collect more data; tune hyper-parameters", $null, "article_9"),
    @("35b3012a-e79a-46bc-b468-2fda279cfe8c", "2025-05-21T18:12:25.589810", 0.65, "DRIFT_VULNERABILITY", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "MEDIUM", "Enable drift monitoring; schedule periodic retraining", $null, "article_9"),
    @("b4421b8e-7d46-4cc2-9b4f-effdb9f3e8a4", "2025-05-21T18:13:46.120087", 0.65, "BIAS_PROTECTED_GROUPS", "Unfair bias against protected demographic groups", "HIGH", "Re-sample training data; add fairness constraints or post-processing techniques", "num__AGE_YEARS: 1.000 disparity`nNAME_EDUCATION_TYPE: 0.410 disparity`nNAME_HOUSING_TYPE: 0.202 disparity`n", "article_10"),
    @("b4421b8e-7d46-4cc2-9b4f-effdb9f3e8a4", "2025-05-21T18:13:46.120087", 0.65, "LOW_ACCURACY", "Model accuracy below 0.75", "MEDIUM", "Collect more data; tune hyper-parameters", $null, "article_9"),
    @("b4421b8e-7d46-4cc2-9b4f-effdb9f3e8a4", "2025-05-21T18:13:46.120087", 0.65, "DRIFT_VULNERABILITY", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "MEDIUM", "Enable drift monitoring; schedule periodic retraining", $null, "article_9"),
    @("ea9382b8-ad1d-41a8-a10b-9b28715e0f77", "2025-05-21T19:06:22.324942", 0.65, "BIAS_PROTECTED_GROUPS", "Unfair bias against protected demographic groups", "HIGH", "Re-sample training data; add fairness constraints or post-processing techniques", "num__AGE_YEARS: 1.000 disparity`nNAME_EDUCATION_TYPE: 0.410 disparity`nNAME_HOUSING_TYPE: 0.202 disparity`n", "article_10"),
    @("ea9382b8-ad1d-41a8-a10b-9b28715e0f77", "2025-05-21T19:06:22.324942", 0.65, "LOW_ACCURACY", "Model accuracy below 0.75", "MEDIUM", "Collect more data; tune hyper-parameters", $null, "article_9"),
    @("ea9382b8-ad1d-41a8-a10b-9b28715e0f77", "2025-05-21T19:06:22.324942", 0.65, "DRIFT_VULNERABILITY", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "MEDIUM", "Enable drift monitoring; schedule periodic retraining", $null, "article_9"),
    @("e244f4ce-832c-4618-addc-984857f50653", "2025-05-21T20:19:38.154359", 0.65, "BIAS_PROTECTED_GROUPS", "Unfair bias against protected demographic groups", "HIGH", "Re-sample training data; add fairness constraints or post-processing techniques", "num__AGE_YEARS: 1.000 disparity`nNAME_EDUCATION_TYPE: 0.410 disparity`nNAME_HOUSING_TYPE: 0.202 disparity`n", "article_10"),
    @("e244f4ce-832c-4618-addc-984857f50653", "2025-05-21T20:19:38.154359", 0.65, "LOW_ACCURACY", "Model accuracy below 0.75", "MEDIUM", "Collect more data; tune hyper-parameters", "", "article_9"),
    @("e244f4ce-832c-4618-addc-984857f50653", "2025-05-21T20:19:38.154359", 0.65, "DRIFT_VULNERABILITY", "ROC-AUC risk proxy > 0.3 indicates drift fragility", "MEDIUM", "Enable drift monitoring; schedule periodic retraining", "", "article_9"),
)

$rowNum = 19
foreach ($row in $hazardsData) {
    $hazards.Cells.Item($rowNum, 1).Value = $row[0]
    $hazards.Cells.Item($rowNum, 2).Value = $row[1]
    $hazards.Cells.Item($rowNum, 3).Value = $row[2]
    $hazards.Cells.Item($rowNum, 4).Value = $row[3]
    $hazards.Cells.Item($rowNum, 5).Value = $row[4]
    $hazards.Cells.Item($rowNum, 6).Value = $row[5]
    $hazards.Cells.Item($rowNum, 7).Value = $row[6]
    # "Details" (H) only has content for bias rows; otherwise left blank
    if ($row[7] -ne $null) {
        $hazards.Cells.Item($rowNum, 8).Value = $row[7]
    }
    $hazards.Cells.Item($rowNum, 9).Value = $row[8]
    $rowNum = $rowNum + 1
}
